# "Generate Report for Handback"
#
# This script fills in the handback results for the zh-cn and de-de
# localization rows: it records the generated target/handback files,
# the handback timestamps, marks the overview/status as handed back,
# and widens a couple of columns so the longer text fits.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

$mdName = "ebe840a9-4e14-4a32-93d4-fe2121ecce10.md"
$mdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/312b8b8d50d511835244b494518324e3296917ca/e2e/ebe840a9-4e14-4a32-93d4-fe2121ecce10.md"

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Status text: flip every "Ready for handoff" cell to the handed-back
# status (Overview!E2/F2, zh-cn!C2, de-de!C2 all shared this text).
# ---------------------------------------------------------------------
$ws_overview.Range("E2").Value = $statusHandedBack
$ws_overview.Range("F2").Value = $statusHandedBack
$ws_zhcn.Range("C2").Value = $statusHandedBack
$ws_dede.Range("C2").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn row: Latest Target File (I2), Latest Handback File (J2) and
# Latest Handback DateTime (K2).
# ---------------------------------------------------------------------
$ws_zhcn.Range("I2").Value = $mdName
$ws_zhcn.Hyperlinks.Add($ws_zhcn.Range("I2"), $mdUrl, "", "", $mdName)
$ws_zhcn.Range("J2").Value = "ebe840a9-4e14-4a32-93d4-fe2121ecce10.a9ed4fc76c8f56ee302467c35ae58749087ab49f.zh-cn.xlf"
$ws_zhcn.Range("K2").Value = "2016-08-18 15:04:32"

# ---------------------------------------------------------------------
# de-de row: Latest Target File (I2), Latest Handback File (J2) and
# Latest Handback DateTime (K2).
# ---------------------------------------------------------------------
$ws_dede.Range("I2").Value = $mdName
$ws_dede.Hyperlinks.Add($ws_dede.Range("I2"), $mdUrl, "", "", $mdName)
$ws_dede.Range("J2").Value = "ebe840a9-4e14-4a32-93d4-fe2121ecce10.a9ed4fc76c8f56ee302467c35ae58749087ab49f.de-de.xlf"
$ws_dede.Range("K2").Value = "2016-08-18 15:04:40"

# ---------------------------------------------------------------------
# Column widths: widen the status columns (now holding the longer
# "Handed back: in sync with en-US" text) and the target/handback file
# columns (now holding the long generated .xlf file names).
# ColumnWidth assignments are snapped by the engine to the nearest
# 1/6 character; the values below are chosen so the saved width comes
# out to the nearest achievable width (30 / 40 characters).
# ---------------------------------------------------------------------
# (numeric column indices are used below -- Columns.Item("<letter>") has
# been observed to mis-marshal through this COM interop and silently
# fail to apply the width)
$ws_overview.Columns.Item(5).ColumnWidth = 29.1666666666667
$ws_overview.Columns.Item(6).ColumnWidth = 29.1666666666667

$ws_zhcn.Columns.Item(3).ColumnWidth = 29.1666666666667
$ws_zhcn.Columns.Item(9).ColumnWidth = 39.1666666666667
$ws_zhcn.Columns.Item(10).ColumnWidth = 39.1666666666667

$ws_dede.Columns.Item(3).ColumnWidth = 29.1666666666667
$ws_dede.Columns.Item(9).ColumnWidth = 39.1666666666667
$ws_dede.Columns.Item(10).ColumnWidth = 39.1666666666667
